$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Pietra Sales -> Felipe Melo
$ws.Cells.Item(2, 1).Value = 38892
$ws.Cells.Item(2, 2).Value = "Felipe Melo"
$ws.Cells.Item(2, 3).Value = "Operacoes"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 45086
$ws.Cells.Item(2, 7).Value = 9819.280000000001

# Row 3: Rael Pereira -> Lunna Pereira
$ws.Cells.Item(3, 1).Value = 73884
$ws.Cells.Item(3, 2).Value = "Lunna Pereira"
$ws.Cells.Item(3, 3).Value = "Juridico"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45089
$ws.Cells.Item(3, 7).Value = 7428.81

# Row 4: Stephany Silva -> Fernando Vasconcelos
$ws.Cells.Item(4, 1).Value = 66692
$ws.Cells.Item(4, 2).Value = "Fernando Vasconcelos"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Consulta medica"
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 45078
$ws.Cells.Item(4, 7).Value = 5016.95

# Row 5: Valentina Costela -> Sr. Vitor Hugo Lima
$ws.Cells.Item(5, 1).Value = 91871
$ws.Cells.Item(5, 2).Value = "Sr. Vitor Hugo Lima"
$ws.Cells.Item(5, 3).Value = "Operacoes"
$ws.Cells.Item(5, 4).Value = "Outros"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 45081
$ws.Cells.Item(5, 7).Value = 7670.13

# Row 6: Helena Novais -> José Pedro Freitas
$ws.Cells.Item(6, 1).Value = 41852
$ws.Cells.Item(6, 2).Value = "José Pedro Freitas"
$ws.Cells.Item(6, 3).Value = "Recursos Humanos"
$ws.Cells.Item(6, 4).Value = "Consulta medica"
$ws.Cells.Item(6, 6).Value = 45102
$ws.Cells.Item(6, 7).Value = 9332.459999999999

# Row 7: Eloá Silva -> Luísa da Rosa
$ws.Cells.Item(7, 1).Value = 33622
$ws.Cells.Item(7, 2).Value = "Luísa da Rosa"
$ws.Cells.Item(7, 3).Value = "P&D"
$ws.Cells.Item(7, 4).Value = "Consulta medica"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 45092
$ws.Cells.Item(7, 7).Value = 4266.76

# Row 8: Sr. Benicio Silva -> Sra. Eloá Almeida
$ws.Cells.Item(8, 1).Value = 42666
$ws.Cells.Item(8, 2).Value = "Sra. Eloá Almeida"
$ws.Cells.Item(8, 4).Value = "Doenca"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45100
$ws.Cells.Item(8, 7).Value = 7321.26

# Row 9: Dr. Cauã Albuquerque -> Luiz Fernando Pereira
$ws.Cells.Item(9, 1).Value = 68850
$ws.Cells.Item(9, 2).Value = "Luiz Fernando Pereira"
$ws.Cells.Item(9, 4).Value = "Viagem de negocios"
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 45085
$ws.Cells.Item(9, 7).Value = 5559.02

# Row 10: Dom da Mota -> Ryan Caldeira
$ws.Cells.Item(10, 1).Value = 18199
$ws.Cells.Item(10, 2).Value = "Ryan Caldeira"
$ws.Cells.Item(10, 3).Value = "Financeiro"
$ws.Cells.Item(10, 4).Value = "Consulta medica"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 45102
$ws.Cells.Item(10, 7).Value = 5535.86

# Row 11: Isis da Paz -> Emilly Cavalcante
$ws.Cells.Item(11, 1).Value = 90172
$ws.Cells.Item(11, 2).Value = "Emilly Cavalcante"
$ws.Cells.Item(11, 3).Value = "P&D"
$ws.Cells.Item(11, 4).Value = "Consulta medica"
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 7).Value = 8739.07
